# Append the new daily profit row (2025-11-15) to the bottom of the table
# on Sheet1, extending the used range from A1:B89 to A1:B90.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 90

# Force column A to be stored as plain text so the "MM/DD/YYYY" string is
# not auto-converted into a date serial number (matching the existing
# rows above it, which are all literal text dates).
$ws.Range("A" + $newRow).NumberFormat = "@"
$ws.Range("A" + $newRow).Value = "11/15/2025"

$ws.Range("B" + $newRow).Value = 8884.26
